$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Find the decorative horizontal-rule paragraphs (a lone <w:r>
#    holding a <w:pict><v:rect .../></w:pict>).  These render as empty
#    paragraphs in the Word object model (style "Normal", Range.Text is
#    just the paragraph mark char(13), with nothing else - the table's
#    own end-of-row empty paragraphs have an extra char(7) so they are
#    not picked up here). There are nine such paragraphs in the
#    document; the last two are left untouched, only the first seven
#    get removed.
#
# 2) Also locate the paragraph beginning with the opening curly quote
#    "Be Good & be Gone." - its style switches from "First Paragraph"
#    to "Body Text".
# ------------------------------------------------------------------

$hrParagraphs = @()
$quoteIndex = -1

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq [char]13) {
        $hrParagraphs += $i
    } elseif ($t.StartsWith([char]8220 + "Be Good")) {
        $quoteIndex = $i
    }
}

# Apply the style change first (index is still valid - no deletions yet).
if ($quoteIndex -gt 0) {
    $d.Paragraphs.Item($quoteIndex).Style = "Body Text"
}

$keepLast = 2
$deleteCount = $hrParagraphs.Count - $keepLast
if ($deleteCount -gt 0) {
    $toDelete = $hrParagraphs[0..($deleteCount - 1)]
    # Delete from the bottom up so earlier indices remain valid.
    for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {
        $idx = $toDelete[$j]
        $d.Paragraphs.Item($idx).Range.Delete()
    }
}
